# Updates odds figures on Sheet1 to match the refreshed FlashScore export.
# Each line below sets one cell's numeric value to the new odd/price.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Ceara vs Ponte Preta)
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 1.85

# Row 3
$ws.Range("R3").Value = 1.53
$ws.Range("V3").Value = 1.57

# Row 4
$ws.Range("R4").Value = 1.36
$ws.Range("V4").Value = 1.5

# Row 8
$ws.Range("S8").Value = 1.4

# Row 9
$ws.Range("G9").Value = 2.35
$ws.Range("I9").Value = 2.75
$ws.Range("J9").Value = 3
$ws.Range("L9").Value = 3.25
$ws.Range("U9").Value = 1.62
$ws.Range("V9").Value = 2.2
$ws.Range("AI9").Value = 10
$ws.Range("AU9").Value = 7

# Row 16
$ws.Range("S16").Value = 1.44
$ws.Range("T16").Value = 2.63

# Row 17
$ws.Range("M17").Value = 1.07
$ws.Range("N17").Value = 9
$ws.Range("O17").Value = 1.36
$ws.Range("P17").Value = 3
$ws.Range("S17").Value = 1.44
$ws.Range("T17").Value = 2.63

# Row 21
$ws.Range("S21").Value = 1.44
$ws.Range("T21").Value = 2.63

# Row 25
$ws.Range("S25").Value = 1.44
$ws.Range("T25").Value = 2.63

# Row 26
$ws.Range("S26").Value = 1.22

# Row 35 (Haverfordwest match block's sibling fixture)
$ws.Range("G35").Value = 1.83
$ws.Range("H35").Value = 3.35
$ws.Range("I35").Value = 4.05
$ws.Range("J35").Value = 2.37
$ws.Range("K35").Value = 2.12
$ws.Range("L35").Value = 4.4
$ws.Range("P35").Value = 2.85
$ws.Range("S35").Value = 1.38
$ws.Range("T35").Value = 2.57
$ws.Range("X35").Value = 8.25
$ws.Range("Z35").Value = 15
$ws.Range("AD35").Value = 6.5
$ws.Range("AG35").Value = 10.75
$ws.Range("AH35").Value = 22
$ws.Range("AI35").Value = 13.5
$ws.Range("AJ35").Value = 65
$ws.Range("AK35").Value = 40
$ws.Range("AN35").Value = 3.65
$ws.Range("AO35").Value = 9
$ws.Range("AP35").Value = 18
$ws.Range("AQ35").Value = 32
$ws.Range("AR35").Value = 65
$ws.Range("AS35").Value = 250
$ws.Range("AT35").Value = 2.55
$ws.Range("AW35").Value = 5.8
$ws.Range("AX35").Value = 23
$ws.Range("BA35").Value = 150
$ws.Range("BB35").Value = 400

# Row 42 (Haverfordwest vs Connahs Q.)
$ws.Range("G42").Value = 2.05
$ws.Range("H42").Value = 3.35
$ws.Range("I42").Value = 3.2
$ws.Range("J42").Value = 2.67
$ws.Range("K42").Value = 2.12
$ws.Range("L42").Value = 3.8
$ws.Range("M42").Value = 1.07
$ws.Range("N42").Value = 6.9
$ws.Range("O42").Value = 1.34
$ws.Range("P42").Value = 3
$ws.Range("Q42").Value = 2.02
$ws.Range("R42").Value = 1.75
$ws.Range("S42").Value = 1.42
$ws.Range("T42").Value = 2.67
$ws.Range("V42").Value = 1.87
$ws.Range("W42").Value = 7
$ws.Range("X42").Value = 9.5
$ws.Range("Y42").Value = 8.75
$ws.Range("Z42").Value = 18.5
$ws.Range("AA42").Value = 17
$ws.Range("AB42").Value = 30
$ws.Range("AC42").Value = 6.9
$ws.Range("AD42").Value = 6.5
$ws.Range("AE42").Value = 15.5
$ws.Range("AG42").Value = 9.25
$ws.Range("AH42").Value = 16.5
$ws.Range("AI42").Value = 11.75
$ws.Range("AJ42").Value = 40
$ws.Range("AK42").Value = 30
$ws.Range("AL42").Value = 40
$ws.Range("AN42").Value = 3.95
$ws.Range("AO42").Value = 10.5
$ws.Range("AP42").Value = 20
$ws.Range("AQ42").Value = 40
$ws.Range("AR42").Value = 80
$ws.Range("AS42").Value = 300
$ws.Range("AT42").Value = 2.67
$ws.Range("AW42").Value = 5.1
$ws.Range("AX42").Value = 18
$ws.Range("AY42").Value = 26
$ws.Range("AZ42").Value = 90
